# add area to Q files stn4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Area column (G), per-segment area, mirrors the Q (E) column pattern ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Totals ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Recalculate + fix up the selection/view to match the edited cells ---
[void]$ws.Range("J2:K2").Select()

$wb.Application.Calculate()
